$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (it sat on an empty paragraph
#    after "... a series of manually set waypoints are adequate.")
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Highlight "(Figure)" in "... fits a Hermite-poly onto them (Figure). From ..."
# ------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Replacement.ClearFormatting()
$rng1.Find.Replacement.Highlight = 1
$rng1.Find.Execute("(Figure)", $false, $false, $false, $false, $false, $true, 1, $false, "(Figure)", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Highlight "(Figure)." in "The calculation of the key-points is detailed in (Figure)."
#    (scope the search to start right after the previous match so it
#    cannot re-match the first occurrence)
# ------------------------------------------------------------------
$rng2 = $d.Range($rng1.End, $d.Content.End)
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$rng2.Find.Replacement.Highlight = 1
$rng2.Find.Execute("(Figure).", $false, $false, $false, $false, $false, $true, 1, $false, "(Figure).", 2) | Out-Null

# Re-create the "_GoBack" bookmark at the very end of that paragraph
# (right after "(Figure)."). A truly zero-length Range can't be
# addressed directly, so append a one-character marker, bookmark that
# character, then delete it again -- this leaves bookmarkStart/End
# collapsed together exactly where the marker was.
$rng2.InsertAfter("X")
$markerRng = $d.Range($rng2.End - 1, $rng2.End)
$d.Bookmarks.Add("_GoBack", $markerRng)
$markerRng = $d.Range($rng2.End - 1, $rng2.End)
$markerRng.Text = ""

# ------------------------------------------------------------------
# 4) Highlight "(Figure)" in ", around the Turning Waypoint. (Figure)"
#    (scope the search to start right after the second match)
# ------------------------------------------------------------------
$rng3 = $d.Range($rng2.End, $d.Content.End)
$rng3.Find.ClearFormatting()
$rng3.Find.Replacement.ClearFormatting()
$rng3.Find.Replacement.Highlight = 1
$rng3.Find.Execute("(Figure)", $false, $false, $false, $false, $false, $true, 1, $false, "(Figure)", 2) | Out-Null
